# Generate Report for Handoff
#
# A new handoff batch was kicked off for the "504921b1-8bf4-445b-8b75-..."
# file (row 7) as well as the two files that depend on it
# ("1f5b0a06-9503-483a-bf30-..." row 10 and "e2700748-ebd0-492f-a585-..."
# row 14). All three rows now share the same "Latest Handoff Datetime"
# (column D) for each localized-language status sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value  = "2016-03-10 20:27:53"
$wsZhCn.Range("D10").Value = "2016-03-10 20:27:53"
$wsZhCn.Range("D14").Value = "2016-03-10 20:27:53"

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value  = "2016-03-10 20:28:00"
$wsDeDe.Range("D10").Value = "2016-03-10 20:28:00"
$wsDeDe.Range("D14").Value = "2016-03-10 20:28:00"
